$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9223046214701632
$ws.Range("J2").Value = 0.9223046214701632
$ws.Range("M2").Value = 1.443038
$ws.Range("N2").Value = 4.329114
$ws.Range("O2").Value = 0.0289666880885598
$ws.Range("P2").Value = 0.0289666880885598
$ws.Range("Q2").Value = 2.407739206798
$ws.Range("R2").Value = 21.669652861182
$ws.Range("S2").Value = 0.02671611029276343
$ws.Range("T2").Value = 0.02671611029276343

# Row 3
$ws.Range("I3").Value = 0.9223046214701632
$ws.Range("J3").Value = 0.9223046214701632
$ws.Range("N3").Value = 87.61054300000001
$ws.Range("O3").Value = 0.5862140087672342
$ws.Range("P3").Value = 0.5862140087672342
$ws.Range("Q3").Value = 48.72667693896767
$ws.Range("R3").Value = 438.5400924507091
$ws.Range("S3").Value = 0.5406678894565708
$ws.Range("T3").Value = 0.5406678894565708

# Row 4
$ws.Range("I4").Value = 0.9223046214701632
$ws.Range("J4").Value = 0.9223046214701632
$ws.Range("M4").Value = 19.170603
$ws.Range("N4").Value = 57.511809
$ws.Range("O4").Value = 0.384819303144206
$ws.Range("P4").Value = 0.384819303144206
$ws.Range("Q4").Value = 31.986553688163
$ws.Range("R4").Value = 287.878983193467
$ws.Range("S4").Value = 0.3549206217208289
$ws.Range("T4").Value = 0.3549206217208289

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.140557
$ws.Range("H5").Value = 0.421671
$ws.Range("I5").Value = 0.07769537852983674
$ws.Range("J5").Value = 0.07769537852983674
$ws.Range("M5").Value = 1.443038
$ws.Range("N5").Value = 4.329114
$ws.Range("O5").Value = 0.0289666880885598
$ws.Range("P5").Value = 0.0289666880885598
$ws.Range("Q5").Value = 0.202829092166
$ws.Range("R5").Value = 1.825461829494
$ws.Range("S5").Value = 0.002250577795796366
$ws.Range("T5").Value = 0.002250577795796366

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.140557
$ws.Range("H6").Value = 0.421671
$ws.Range("I6").Value = 0.07769537852983674
$ws.Range("J6").Value = 0.07769537852983674
$ws.Range("N6").Value = 87.61054300000001
$ws.Range("O6").Value = 0.5862140087672342
$ws.Range("P6").Value = 0.5862140087672342
$ws.Range("Q6").Value = 4.104758364150334
$ws.Range("R6").Value = 36.942825277353
$ws.Range("S6").Value = 0.04554611931066329
$ws.Range("T6").Value = 0.04554611931066329

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.140557
$ws.Range("H7").Value = 0.421671
$ws.Range("I7").Value = 0.07769537852983674
$ws.Range("J7").Value = 0.07769537852983674
$ws.Range("M7").Value = 19.170603
$ws.Range("N7").Value = 57.511809
$ws.Range("O7").Value = 0.384819303144206
$ws.Range("P7").Value = 0.384819303144206
$ws.Range("Q7").Value = 2.694562445871
$ws.Range("R7").Value = 24.251062012839
$ws.Range("S7").Value = 0.02989868142337708
$ws.Range("T7").Value = 0.02989868142337708
